$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("extraction consolidation result")

# Update reference IDs after adding new papers: S26 -> S27
$ws.Range("B3").Value = "S27"
$ws.Range("B4").Value = "S27"

$ws.Range("B5").Select() | Out-Null
